$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create a new number format for date/time display: YYYY-MM-DD HH:MM:SS
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Years corresponding to rows 2..39 (1987Q4 .. 2024Q4), each mapped to Dec 31 of that year
$startYear = 1987
$row = 2
for ($year = $startYear; $year -le 2024; $year++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = $dateFormat
    $cell.Value = (Get-Date -Year $year -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0).Date
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $row++
}
